$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps its text (string) type even when new values look numeric
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '62.332.11'
$ws.Range("E2").Value = '  -1.67%  '
$ws.Range("D3").Value = '3.020.50'
$ws.Range("E3").Value = '  -1.48%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '584.83'
$ws.Range("E5").Value = '  -0.94%  '
$ws.Range("D6").Value = '147.51'
$ws.Range("E6").Value = '  -4.01%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '0.525'
$ws.Range("E8").Value = '  -2.38%  '
$ws.Range("D9").Value = '3.015.62'
$ws.Range("E9").Value = '  -1.55%  '
$ws.Range("D10").Value = '0.150'
$ws.Range("E10").Value = '  -4.11%  '
$ws.Range("D11").Value = '5.83'
$ws.Range("E11").Value = '  -1.30%  '
$ws.Range("D12").Value = '0.444'
$ws.Range("E12").Value = '  -1.42%  '
$ws.Range("D13").Value = '0.0000229'
$ws.Range("E13").Value = '  -3.49%  '
$ws.Range("D14").Value = '34.62'
$ws.Range("E14").Value = '  -5.16%  '
$ws.Range("E15").Value = '  +1.96%  '
$ws.Range("D16").Value = '3.519.07'
$ws.Range("E16").Value = '  -1.48%  '
$ws.Range("D17").Value = '7.08'
$ws.Range("E17").Value = '  -1.43%  '
$ws.Range("D18").Value = '62.303.84'
$ws.Range("E18").Value = '  -1.53%  '
$ws.Range("D19").Value = '3.020.83'
$ws.Range("E19").Value = '  -1.30%  '
$ws.Range("D20").Value = '463.54'
$ws.Range("E20").Value = '  -4.26%  '
$ws.Range("D21").Value = '13.99'
$ws.Range("E21").Value = '  -3.41%  '
$ws.Range("D22").Value = '0.685'
$ws.Range("E23").Value = '  -0.78%  '
$ws.Range("D24").Value = '81.76'
$ws.Range("E24").Value = '  -0.38%  '
$ws.Range("D25").Value = '2.26'
$ws.Range("E25").Value = '  -5.55%  '
$ws.Range("D26").Value = '12.30'
$ws.Range("E26").Value = '  -3.82%  '
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("D28").Value = '9.99'
$ws.Range("E28").Value = '  -5.05%  '
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.26%  '
$ws.Range("D30").Value = '2.63'
$ws.Range("E30").Value = '  -2.13%  '
$ws.Range("E31").Value = '  -4.90%  '
$ws.Range("D32").Value = '29.19'
$ws.Range("E32").Value = '  +7.15%  '
$ws.Range("D33").Value = '2.10'
$ws.Range("E33").Value = '  -5.64%  '
$ws.Range("D34").Value = '0.108'
$ws.Range("E34").Value = '  -2.65%  '
$ws.Range("B35").Value = 'PEPE'
$ws.Range("C35").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D35").Value = '0.0₃0798'
$ws.Range("E35").Value = '  -2.61%  '
$ws.Range("B36").Value = 'Mantle'
$ws.Range("C36").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D36").Value = '1.03'
$ws.Range("E36").Value = '  -3.08%  '
$ws.Range("E37").Value = '  -4.49%  '
$ws.Range("D38").Value = '2.12'
$ws.Range("E38").Value = '  -4.24%  '
$ws.Range("D39").Value = '50.35'
$ws.Range("E39").Value = '  -0.57%  '
$ws.Range("D40").Value = '9.05'
$ws.Range("E40").Value = '  -2.69%  '
$ws.Range("D41").Value = '2.95'
$ws.Range("E41").Value = '  -9.22%  '
$ws.Range("E42").Value = '  -0.41%  '
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").Value = '388.77'
$ws.Range("E43").Value = '  -11.69%  '
$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").Value = '0.275'
$ws.Range("E44").Value = '  -4.79%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").Value = '0.0357'
$ws.Range("E45").Value = '  -1.48%  '
$ws.Range("D46").Value = '2.756.51'
$ws.Range("E46").Value = '  -2.24%  '
$ws.Range("D47").Value = '37.49'
$ws.Range("E47").Value = '  -5.95%  '
$ws.Range("D48").Value = '128.29'
$ws.Range("E48").Value = '  -3.30%  '
$ws.Range("E50").Value = '  -0.99%  '
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").Value = '24.08'
$ws.Range("E51").Value = '  -4.89%  '

# Restore default (Normal) style so no stray number-format styling remains
$ws.Range("D2:D51").Style = "Normal"
